$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44435
$ws.Range("I2").Value2 = 'Primera'
$ws.Range("J2").Value2 = 25
$ws.Range("K2").Value2 = 14000
$ws.Range("L2").Value2 = 14000
$ws.Range("M2").Value2 = 14000
$ws.Range("N2").Value2 = '$/caja 30 unidades'
$ws.Range("O2").Value2 = 'Provincia de Limarí'
$ws.Range("P2").Value2 = 467
$ws.Range("Q2").Value2 = 30

# Row 3
$ws.Range("D3").Value2 = 44435
$ws.Range("I3").Value2 = 'Primera'
$ws.Range("J3").Value2 = 25
$ws.Range("K3").Value2 = 14000
$ws.Range("L3").Value2 = 14000
$ws.Range("M3").Value2 = 14000
$ws.Range("N3").Value2 = '$/caja 30 unidades'
$ws.Range("O3").Value2 = 'Provincia del Elquí'
$ws.Range("P3").Value2 = 467
$ws.Range("Q3").Value2 = 30

# Row 4
$ws.Range("D4").Value2 = 44432
$ws.Range("I4").Value2 = 'Primera'
$ws.Range("J4").Value2 = 25
$ws.Range("K4").Value2 = 14000
$ws.Range("L4").Value2 = 14000
$ws.Range("M4").Value2 = 14000
$ws.Range("N4").Value2 = '$/caja 30 unidades'
$ws.Range("O4").Value2 = 'Provincia del Elquí'
$ws.Range("P4").Value2 = 467
$ws.Range("Q4").Value2 = 30

# Row 5
$ws.Range("D5").Value2 = 44467
$ws.Range("I5").Value2 = 'Primera'
$ws.Range("J5").Value2 = 35
$ws.Range("K5").Value2 = 12000
$ws.Range("L5").Value2 = 12000
$ws.Range("M5").Value2 = 12000
$ws.Range("N5").Value2 = '$/caja 30 unidades'
$ws.Range("O5").Value2 = 'Provincia de Limarí'
$ws.Range("P5").Value2 = 400
$ws.Range("Q5").Value2 = 30

# Row 6
$ws.Range("D6").Value2 = 44418
$ws.Range("I6").Value2 = 'Primera'
$ws.Range("J6").Value2 = 30
$ws.Range("K6").Value2 = 15000
$ws.Range("L6").Value2 = 15000
$ws.Range("M6").Value2 = 15000
$ws.Range("N6").Value2 = '$/caja 30 unidades'
$ws.Range("O6").Value2 = 'Provincia de Limarí'
$ws.Range("P6").Value2 = 500
$ws.Range("Q6").Value2 = 30

# Row 7
$ws.Range("D7").Value2 = 44474
$ws.Range("I7").Value2 = 'Primera'
$ws.Range("J7").Value2 = 45
$ws.Range("K7").Value2 = 10000
$ws.Range("L7").Value2 = 10000
$ws.Range("M7").Value2 = 10000
$ws.Range("N7").Value2 = '$/caja 30 unidades'
$ws.Range("O7").Value2 = 'Provincia de Limarí'
$ws.Range("P7").Value2 = 333
$ws.Range("Q7").Value2 = 30

# Row 8
$ws.Range("D8").Value2 = 44446
$ws.Range("I8").Value2 = 'Primera'
$ws.Range("J8").Value2 = 25
$ws.Range("K8").Value2 = 14000
$ws.Range("L8").Value2 = 14000
$ws.Range("M8").Value2 = 14000
$ws.Range("N8").Value2 = '$/caja 30 unidades'
$ws.Range("O8").Value2 = 'Provincia de Limarí'
$ws.Range("P8").Value2 = 467
$ws.Range("Q8").Value2 = 30

# Row 9
$ws.Range("D9").Value2 = 44460
$ws.Range("I9").Value2 = 'Primera'
$ws.Range("J9").Value2 = 45
$ws.Range("K9").Value2 = 13000
$ws.Range("L9").Value2 = 13000
$ws.Range("M9").Value2 = 13000
$ws.Range("N9").Value2 = '$/caja 30 unidades'
$ws.Range("O9").Value2 = 'Provincia de Limarí'
$ws.Range("P9").Value2 = 433
$ws.Range("Q9").Value2 = 30

# Row 10
$ws.Range("D10").Value2 = 44376
$ws.Range("I10").Value2 = 'Primera'
$ws.Range("J10").Value2 = 25
$ws.Range("K10").Value2 = 18000
$ws.Range("L10").Value2 = 18000
$ws.Range("M10").Value2 = 18000
$ws.Range("N10").Value2 = '$/caja 30 unidades'
$ws.Range("O10").Value2 = 'Provincia de Limarí'
$ws.Range("P10").Value2 = 600
$ws.Range("Q10").Value2 = 30

# Row 11
$ws.Range("D11").Value2 = 44453
$ws.Range("I11").Value2 = 'Primera'
$ws.Range("J11").Value2 = 50
$ws.Range("K11").Value2 = 12000
$ws.Range("L11").Value2 = 12000
$ws.Range("M11").Value2 = 12000
$ws.Range("N11").Value2 = '$/caja 30 unidades'
$ws.Range("O11").Value2 = 'Provincia de Limarí'
$ws.Range("P11").Value2 = 400
$ws.Range("Q11").Value2 = 30

# Row 12
$ws.Range("D12").Value2 = 44421
$ws.Range("I12").Value2 = 'Primera'
$ws.Range("J12").Value2 = 25
$ws.Range("K12").Value2 = 15000
$ws.Range("L12").Value2 = 16000
$ws.Range("M12").Value2 = 15400
$ws.Range("N12").Value2 = '$/caja 30 unidades'
$ws.Range("O12").Value2 = 'Provincia de Limarí'
$ws.Range("P12").Value2 = 513
$ws.Range("Q12").Value2 = 30

# Row 13
$ws.Range("D13").Value2 = 44841
$ws.Range("I13").Value2 = 'Primera'
$ws.Range("J13").Value2 = 45
$ws.Range("K13").Value2 = 12000
$ws.Range("L13").Value2 = 12000
$ws.Range("M13").Value2 = 12000
$ws.Range("N13").Value2 = '$/caja 30 unidades'
$ws.Range("O13").Value2 = 'Provincia de Limarí'
$ws.Range("P13").Value2 = 400
$ws.Range("Q13").Value2 = 30

# Row 14
$ws.Range("D14").Value2 = 44841
$ws.Range("I14").Value2 = 'Segunda'
$ws.Range("J14").Value2 = 45
$ws.Range("K14").Value2 = 10000
$ws.Range("L14").Value2 = 10000
$ws.Range("M14").Value2 = 10000
$ws.Range("N14").Value2 = '$/caja 40 unidades'
$ws.Range("O14").Value2 = 'Provincia de Limarí'
$ws.Range("P14").Value2 = 250
$ws.Range("Q14").Value2 = 40

# Row 15
$ws.Range("D15").Value2 = 44449
$ws.Range("I15").Value2 = 'Primera'
$ws.Range("J15").Value2 = 45
$ws.Range("K15").Value2 = 12000
$ws.Range("L15").Value2 = 12000
$ws.Range("M15").Value2 = 12000
$ws.Range("N15").Value2 = '$/caja 30 unidades'
$ws.Range("O15").Value2 = 'Provincia de Limarí'
$ws.Range("P15").Value2 = 400
$ws.Range("Q15").Value2 = 30

# Row 16
$ws.Range("D16").Value2 = 44425
$ws.Range("I16").Value2 = 'Primera'
$ws.Range("J16").Value2 = 35
$ws.Range("K16").Value2 = 14000
$ws.Range("L16").Value2 = 14000
$ws.Range("M16").Value2 = 14000
$ws.Range("N16").Value2 = '$/caja 30 unidades'
$ws.Range("O16").Value2 = 'Provincia de Limarí'
$ws.Range("P16").Value2 = 467
$ws.Range("Q16").Value2 = 30
